$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 282. This shifts all the
# existing data rows (282..309) down to (283..310) and extends the
# sheet dimension accordingly, matching the "A1:R309" -> "A1:R310" change.
$ws.Rows("282:282").Insert()

# Populate the freshly inserted row 282 with the new weekly record.
$ws.Range("A282").Value = 9
$ws.Range("B282").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C282").Value = "Metropolitana"
$ws.Range("D282").Value = 44826
$ws.Range("E282").Value = 13
$ws.Range("F282").Value = 100112001
$ws.Range("G282").Value = "Berenjena"
$ws.Range("H282").Value = "Sin especificar"
$ws.Range("I282").Value = "Primera"
$ws.Range("J282").Value = 90
$ws.Range("K282").Value = 12000
$ws.Range("L282").Value = 12000
$ws.Range("M282").Value = 12000
$ws.Range("N282").Value = "`$/caja 50 unidades"
$ws.Range("O282").Value = "Región de Arica y Parinacota"
$ws.Range("P282").Value = 240
$ws.Range("Q282").Value = 50
$ws.Range("R282").Value = "Hortaliza"
